# Update the aircraft table: row 7's model is renamed RF12 -> RF32 and its
# passenger capacity / full-tanks values are both updated to 600.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "RF32"
$ws.Range("B7").Value = 600
$ws.Range("C7").Value = 600
